# Update sheet data for the week ending 2022-03-03 (carjacking by month, YoY)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the "through" label to reflect the new date
$ws.Name = "Through 2022-02-23"
$ws.Range("A3").Value = "February (through 02-23)"

# Update February row (row 3) values for years 2017-2022 (columns D-I)
$ws.Range("D3").Value = 47
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 25
$ws.Range("G3").Value = 59
$ws.Range("H3").Value = 103
$ws.Range("I3").Value = 118

# Update Total row (row 4) values for years 2017-2022 (columns D-I)
$ws.Range("D4").Value = 122
$ws.Range("E4").Value = 133
$ws.Range("F4").Value = 74
$ws.Range("G4").Value = 133
$ws.Range("H4").Value = 320
$ws.Range("I4").Value = 277
